# Update NOAA-derived average_county_temperature (K) values and the downstream
# worst/best ASHP COP columns (R, S) that were recomputed from the new temperatures.
# NAICS 311230 rows (re-added facility) and other merged-dataset rows are included below.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: K2=13.0158303464755, R2=1.853964204859962, S2=2.02833814451736
$ws.Range("K2").Value = 13.0158303464755
$ws.Range("R2").Value = 1.853964204859962
$ws.Range("S2").Value = 2.02833814451736

# Row 3: K3=13.0158303464755, R3=3.15624931769735, S3=3.848474349579903
$ws.Range("K3").Value = 13.0158303464755
$ws.Range("R3").Value = 3.15624931769735
$ws.Range("S3").Value = 3.848474349579903

# Row 6: K6=19.79629629629628, R6=1.98600466835246, S6=2.18975222777657
$ws.Range("K6").Value = 19.79629629629628
$ws.Range("R6").Value = 1.98600466835246
$ws.Range("S6").Value = 2.18975222777657

# Row 7: K7=19.79629629629628, R7=3.629680458828347, S7=4.589715938979482
$ws.Range("K7").Value = 19.79629629629628
$ws.Range("R7").Value = 3.629680458828347
$ws.Range("S7").Value = 4.589715938979482

# Row 15: K15=21.79166666666666, R15=2.028520339740724, S15=2.242263395092639
$ws.Range("K15").Value = 21.79166666666666
$ws.Range("R15").Value = 2.028520339740724
$ws.Range("S15").Value = 2.242263395092639

# Row 16: K16=21.79166666666666, R16=3.797299903567984, S16=4.865495608531995
$ws.Range("K16").Value = 21.79166666666666
$ws.Range("R16").Value = 3.797299903567984
$ws.Range("S16").Value = 4.865495608531995

# Row 17: K17=21.79166666666666
$ws.Range("K17").Value = 21.79166666666666

# Row 20: K20=15.36574074074072, R20=1.897690627758933, S20=2.081514384587572
$ws.Range("K20").Value = 15.36574074074072
$ws.Range("R20").Value = 1.897690627758933
$ws.Range("S20").Value = 2.081514384587572

# Row 21: K21=15.36574074074072, R21=3.305680440257437, S21=4.076649924074289
$ws.Range("K21").Value = 15.36574074074072
$ws.Range("R21").Value = 3.305680440257437
$ws.Range("S21").Value = 4.076649924074289

# Row 28: K28=12.51681286549706, R28=3.126239257907711, S28=3.803269267167952
$ws.Range("K28").Value = 12.51681286549706
$ws.Range("R28").Value = 3.126239257907711
$ws.Range("S28").Value = 3.803269267167952

# Row 29: K29=12.51681286549706, R29=1.844936767548521, S29=2.017393709936214
$ws.Range("K29").Value = 12.51681286549706
$ws.Range("R29").Value = 1.844936767548521
$ws.Range("S29").Value = 2.017393709936214

# Row 30: K30=19.79629629629628, R30=1.98600466835246, S30=2.18975222777657
$ws.Range("K30").Value = 19.79629629629628
$ws.Range("R30").Value = 1.98600466835246
$ws.Range("S30").Value = 2.18975222777657

# Row 31: K31=19.79629629629628, R31=3.629680458828347, S31=4.589715938979482
$ws.Range("K31").Value = 19.79629629629628
$ws.Range("R31").Value = 3.629680458828347
$ws.Range("S31").Value = 4.589715938979482

# Row 40: K40=21.19907407407406, R40=2.015705049109126, S40=2.22640738080769
$ws.Range("K40").Value = 21.19907407407406
$ws.Range("R40").Value = 2.015705049109126
$ws.Range("S40").Value = 2.22640738080769

# Row 41: K41=21.19907407407406, R41=3.745925377867032, S41=4.780194493904943
$ws.Range("K41").Value = 21.19907407407406
$ws.Range("R41").Value = 3.745925377867032
$ws.Range("S41").Value = 4.780194493904943

# Row 42: K42=21.19907407407406, R42=3.745925377867032, S42=4.780194493904943
$ws.Range("K42").Value = 21.19907407407406
$ws.Range("R42").Value = 3.745925377867032
$ws.Range("S42").Value = 4.780194493904943

# Row 43: K43=21.19907407407406
$ws.Range("K43").Value = 21.19907407407406

# Row 44: K44=21.19907407407406, R44=2.015705049109126, S44=2.22640738080769
$ws.Range("K44").Value = 21.19907407407406
$ws.Range("R44").Value = 2.015705049109126
$ws.Range("S44").Value = 2.22640738080769

# Row 45: K45=21.19907407407406, R45=3.745925377867032, S45=4.780194493904943
$ws.Range("K45").Value = 21.19907407407406
$ws.Range("R45").Value = 3.745925377867032
$ws.Range("S45").Value = 4.780194493904943

# Row 46: K46=21.19907407407406, R46=2.015705049109126, S46=2.22640738080769
$ws.Range("K46").Value = 21.19907407407406
$ws.Range("R46").Value = 2.015705049109126
$ws.Range("S46").Value = 2.22640738080769

# Row 51: K51=21.19907407407406, R51=2.015705049109126, S51=2.22640738080769
$ws.Range("K51").Value = 21.19907407407406
$ws.Range("R51").Value = 2.015705049109126
$ws.Range("S51").Value = 2.22640738080769

# Row 52: K52=21.19907407407406, R52=3.745925377867032, S52=4.780194493904943
$ws.Range("K52").Value = 21.19907407407406
$ws.Range("R52").Value = 3.745925377867032
$ws.Range("S52").Value = 4.780194493904943

# Row 68: K68=15.74228395061728, R68=3.33095021773865, S68=4.115751405322535
$ws.Range("K68").Value = 15.74228395061728
$ws.Range("R68").Value = 3.33095021773865
$ws.Range("S68").Value = 4.115751405322535

# Row 69: K69=15.74228395061728, R69=1.904889690449167, S69=2.090295475371289
$ws.Range("K69").Value = 15.74228395061728
$ws.Range("R69").Value = 1.904889690449167
$ws.Range("S69").Value = 2.090295475371289

# Row 79: K79=15.36574074074072, R79=3.305680440257437, S79=4.076649924074289
$ws.Range("K79").Value = 15.36574074074072
$ws.Range("R79").Value = 3.305680440257437
$ws.Range("S79").Value = 4.076649924074289

# Row 80: K80=15.36574074074072, R80=1.897690627758933, S80=2.081514384587572
$ws.Range("K80").Value = 15.36574074074072
$ws.Range("R80").Value = 1.897690627758933
$ws.Range("S80").Value = 2.081514384587572

# Row 81: K81=3.38888888888889, R81=2.66307484220018, S81=3.130624327233584
$ws.Range("K81").Value = 3.38888888888889
$ws.Range("R81").Value = 2.66307484220018
$ws.Range("S81").Value = 3.130624327233584

# Row 82: K82=3.38888888888889, R82=1.694051767048283, S82=1.836167304537999
$ws.Range("K82").Value = 3.38888888888889
$ws.Range("R82").Value = 1.694051767048283
$ws.Range("S82").Value = 1.836167304537999

# Row 83: K83=2.356481481481501, R83=2.619185573867416, S83=3.069228739776626
$ws.Range("K83").Value = 2.356481481481501
$ws.Range("R83").Value = 2.619185573867416
$ws.Range("S83").Value = 3.069228739776626

# Row 84: K84=2.356481481481501, R84=1.678525338046114, S84=1.817698795724144
$ws.Range("K84").Value = 2.356481481481501
$ws.Range("R84").Value = 1.678525338046114
$ws.Range("S84").Value = 1.817698795724144

# Row 100: K100=13.0158303464755, R100=1.853964204859962, S100=2.02833814451736
$ws.Range("K100").Value = 13.0158303464755
$ws.Range("R100").Value = 1.853964204859962
$ws.Range("S100").Value = 2.02833814451736

# Row 101: K101=13.0158303464755, R101=3.15624931769735, S101=3.848474349579903
$ws.Range("K101").Value = 13.0158303464755
$ws.Range("R101").Value = 3.15624931769735
$ws.Range("S101").Value = 3.848474349579903

# Row 102: K102=18.89814814814816, R102=1.967443877059447, S102=2.16691042047532
$ws.Range("K102").Value = 18.89814814814816
$ws.Range("R102").Value = 1.967443877059447
$ws.Range("S102").Value = 2.16691042047532

# Row 103: K103=18.89814814814816, R103=3.558967664189598, S103=4.475532187740448
$ws.Range("K103").Value = 18.89814814814816
$ws.Range("R103").Value = 3.558967664189598
$ws.Range("S103").Value = 4.475532187740448

# Row 104: K104=19.30324074074072
$ws.Range("K104").Value = 19.30324074074072

# Row 105: K105=19.30324074074072
$ws.Range("K105").Value = 19.30324074074072

# Row 106: K106=19.30324074074072
$ws.Range("K106").Value = 19.30324074074072

# Row 107: K107=19.30324074074072, R107=1.975772235794973, S107=2.177153507468733
$ws.Range("K107").Value = 19.30324074074072
$ws.Range("R107").Value = 1.975772235794973
$ws.Range("S107").Value = 2.177153507468733

# Row 108: K108=19.30324074074072, R108=3.590517197710347, S108=4.526321250243172
$ws.Range("K108").Value = 19.30324074074072
$ws.Range("R108").Value = 3.590517197710347
$ws.Range("S108").Value = 4.526321250243172

# Row 113: K113=13.75752314814816, R113=1.867546171126113, S113=2.044826120875009
$ws.Range("K113").Value = 13.75752314814816
$ws.Range("R113").Value = 1.867546171126113
$ws.Range("S113").Value = 2.044826120875009

# Row 114: K114=13.75752314814816, R114=3.201933436480062, S114=3.917684201664166
$ws.Range("K114").Value = 13.75752314814816
$ws.Range("R114").Value = 3.201933436480062
$ws.Range("S114").Value = 3.917684201664166

# Row 115: K115=13.76976495726495, R115=1.867772014163364, S115=2.045100507661769
$ws.Range("K115").Value = 13.76976495726495
$ws.Range("R115").Value = 1.867772014163364
$ws.Range("S115").Value = 2.045100507661769

# Row 116: K116=13.76976495726495, R116=3.202698560003336, S116=3.918847414586112
$ws.Range("K116").Value = 13.76976495726495
$ws.Range("R116").Value = 3.202698560003336
$ws.Range("S116").Value = 3.918847414586112

# Row 121: K121=-1.819444444444444, R121=1.618523362263702, S121=1.746638928617865
$ws.Range("K121").Value = -1.819444444444444
$ws.Range("R121").Value = 1.618523362263702
$ws.Range("S121").Value = 1.746638928617865

# Row 122: K122=-1.819444444444444, R122=2.455497817501559, S122=2.843656807626497
$ws.Range("K122").Value = -1.819444444444444
$ws.Range("R122").Value = 2.455497817501559
$ws.Range("S122").Value = 2.843656807626497

# Row 129: K129=21.19907407407406, R129=2.015705049109126, S129=2.22640738080769
$ws.Range("K129").Value = 21.19907407407406
$ws.Range("R129").Value = 2.015705049109126
$ws.Range("S129").Value = 2.22640738080769

# Row 130: K130=21.19907407407406, R130=3.745925377867032, S130=4.780194493904943
$ws.Range("K130").Value = 21.19907407407406
$ws.Range("R130").Value = 3.745925377867032
$ws.Range("S130").Value = 4.780194493904943
